$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "carol menezes"
$ws.Range("B3").Value = "111.000.111-00"
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = "F"
$ws.Range("E3").Value = "vila velha"
